$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - force text format to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.222.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.772.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.769.27"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.406.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.776.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.251.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000143"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.920.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.727.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "404.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.60"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  -4.84%  "
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("E39").Value = "  -6.89%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("E51").Value = "  +0.44%  "
